$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 3950.1667
$ws.Range("I74").Value = 3914.5715
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3914.5715
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -2978.5715
$ws.Range("N74").Value = -5872
# Row 77
$ws.Range("H77").Value = 3950.1667
$ws.Range("I77").Value = 3914.5715
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 19572.8575
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -14892.8575
$ws.Range("N77").Value = -29360
# Row 113
$ws.Range("H113").Value = 3770
$ws.Range("I113").Value = 3850.8333
$ws.Range("K113").Value = 3850.8333
$ws.Range("M113").Value = -596.8332999999998
# Row 137
$ws.Range("H137").Value = 1962792.9
$ws.Range("I137").Value = 2942818.5
$ws.Range("K137").Value = 8828455.5
$ws.Range("M137").Value = -8825905.5
# Row 138
$ws.Range("H138").Value = 2827183.8
$ws.Range("I138").Value = 3750
$ws.Range("J138").Value = 3032524.5
$ws.Range("K138").Value = 11250
$ws.Range("L138").Value = 9097573.5
$ws.Range("M138").Value = -6110
$ws.Range("N138").Value = -9107853.5

$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 1383.5
$ws.Range("I37").Value = 1383.5
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1383.5
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -1110.5
# Row 61
$ws.Range("H61").Value = 100200640
$ws.Range("I61").Value = 125125520
$ws.Range("K61").Value = 125125520
$ws.Range("M61").Value = -125125308
# Row 63
$ws.Range("H63").Value = 3368.2778
$ws.Range("I63").Value = 2329.9092
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 2329.9092
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -1643.9092
$ws.Range("N63").Value = -6372
# Row 66
$ws.Range("H66").Value = 3368.2778
$ws.Range("I66").Value = 2329.9092
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 11649.546
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -8217.546
$ws.Range("N66").Value = -31864
# Row 80
$ws.Range("H80").Value = 40699.855
$ws.Range("I80").Value = 28333.334
$ws.Range("J80").Value = 49974.75
$ws.Range("K80").Value = 28333.334
$ws.Range("L80").Value = 49974.75
$ws.Range("M80").Value = -27335.334
$ws.Range("N80").Value = -51970.75
# Row 83
$ws.Range("H83").Value = 40699.855
$ws.Range("I83").Value = 28333.334
$ws.Range("J83").Value = 49974.75
$ws.Range("K83").Value = 85000.00199999999
$ws.Range("L83").Value = 149924.25
$ws.Range("M83").Value = -80008.00199999999
$ws.Range("N83").Value = -159908.25
# Row 136
$ws.Range("H136").Value = 100200640
$ws.Range("I136").Value = 125125520
$ws.Range("K136").Value = 375376560
$ws.Range("M136").Value = -375374010

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1207
$ws.Range("I20").Value = 992.25
$ws.Range("J20").Value = 1314.375
$ws.Range("K20").Value = 992.25
$ws.Range("L20").Value = 1314.375
$ws.Range("M20").Value = -745.25
$ws.Range("N20").Value = -1808.375
# Row 35
$ws.Range("H35").Value = 1900
$ws.Range("I35").Value = 1900
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1900
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -1590
# Row 80
$ws.Range("H80").Value = 579.4545000000001
$ws.Range("I80").Value = 200
$ws.Range("J80").Value = 617.4
$ws.Range("K80").Value = 200
$ws.Range("L80").Value = 617.4
$ws.Range("M80").Value = 798
$ws.Range("N80").Value = -2613.4
# Row 82
$ws.Range("H82").Value = 21333
$ws.Range("I82").Value = 11999.5
$ws.Range("K82").Value = 11999.5
$ws.Range("M82").Value = -11616.5
# Row 83
$ws.Range("H83").Value = 579.4545000000001
$ws.Range("I83").Value = 200
$ws.Range("J83").Value = 617.4
$ws.Range("K83").Value = 1000
$ws.Range("L83").Value = 3087
$ws.Range("M83").Value = 3992
$ws.Range("N83").Value = -13071
# Row 85
$ws.Range("H85").Value = 21333
$ws.Range("I85").Value = 11999.5
$ws.Range("K85").Value = 11999.5
$ws.Range("M85").Value = -10673.5
# Row 99
$ws.Range("H99").Value = 1025.2727
$ws.Range("I99").Value = 952.4286
$ws.Range("J99").Value = 1152.75
$ws.Range("K99").Value = 952.4286
$ws.Range("L99").Value = 1152.75
$ws.Range("M99").Value = 545.5714
$ws.Range("N99").Value = -4148.75
# Row 124
$ws.Range("H124").Value = 37850
$ws.Range("J124").Value = 37850
$ws.Range("L124").Value = 37850
$ws.Range("N124").Value = -47670
# Row 134
$ws.Range("H134").Value = 2124.842
$ws.Range("I134").Value = 2219.111
$ws.Range("K134").Value = 6657.333
$ws.Range("M134").Value = -4122.333
# Row 140
$ws.Range("H140").Value = 59975.26
$ws.Range("J140").Value = 59975.26
$ws.Range("L140").Value = 59975.26
$ws.Range("N140").Value = -70335.26000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 92
$ws.Range("H92").Value = 36638.75
$ws.Range("J92").Value = 36638.75
$ws.Range("L92").Value = 36638.75
$ws.Range("N92").Value = -41630.75
# Row 107
$ws.Range("H107").Value = 1070
$ws.Range("I107").Value = 600
$ws.Range("K107").Value = 600
$ws.Range("M107").Value = 1320

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 24471.643
$ws.Range("I5").Value = 31712.625
$ws.Range("J5").Value = 1300.5
$ws.Range("K5").Value = 95137.875
$ws.Range("L5").Value = 3901.5
$ws.Range("M5").Value = -95025.875
$ws.Range("N5").Value = -4125.5
# Row 107
$ws.Range("H107").Value = 835.14703
$ws.Range("I107").Value = 652.3684
$ws.Range("J107").Value = 1066.6666
$ws.Range("K107").Value = 1957.1052
$ws.Range("L107").Value = 3199.9998
$ws.Range("M107").Value = -37.10519999999997
$ws.Range("N107").Value = -7039.9998
# Row 115
$ws.Range("H115").Value = 2902.25
$ws.Range("I115").Value = 2304.5
$ws.Range("K115").Value = 6913.5
$ws.Range("M115").Value = -5738.5
# Row 135
$ws.Range("H135").Value = 24471.643
$ws.Range("I135").Value = 31712.625
$ws.Range("J135").Value = 1300.5
$ws.Range("K135").Value = 285413.625
$ws.Range("L135").Value = 11704.5
$ws.Range("M135").Value = -282878.625
$ws.Range("N135").Value = -16774.5

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1955.8334
$ws.Range("I113").Value = 1828.6666
$ws.Range("J113").Value = 2083
$ws.Range("K113").Value = 1828.6666
$ws.Range("L113").Value = 2083
$ws.Range("M113").Value = 341.3334
$ws.Range("N113").Value = -6423
# Row 123
$ws.Range("H123").Value = 23552.334
$ws.Range("J123").Value = 23552.334
$ws.Range("L123").Value = 23552.334
$ws.Range("N123").Value = -28452.334

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2357.7693
$ws.Range("I40").Value = 2241
$ws.Range("K40").Value = 2241
$ws.Range("M40").Value = -2105
# Row 61
$ws.Range("H61").Value = 3966.6667
$ws.Range("I61").Value = 3966.6667
$ws.Range("K61").Value = 3966.6667
$ws.Range("M61").Value = -3764.6667
# Row 68
$ws.Range("H68").Value = 1635.9445
$ws.Range("J68").Value = 1850
$ws.Range("L68").Value = 1850
$ws.Range("N68").Value = -3348
# Row 71
$ws.Range("H71").Value = 1635.9445
$ws.Range("J71").Value = 1850
$ws.Range("L71").Value = 9250
$ws.Range("N71").Value = -16738
# Row 113
$ws.Range("H113").Value = 3966.6667
$ws.Range("I113").Value = 3966.6667
$ws.Range("K113").Value = 3966.6667
$ws.Range("M113").Value = -1796.6667
# Row 139
$ws.Range("H139").Value = 43715
$ws.Range("J139").Value = 43715
$ws.Range("L139").Value = 43715
$ws.Range("N139").Value = -53995
